$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1. Rename existing sheets
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Spratelloides"

$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "Gerres"

# ---------------------------------------------------------------
# 2. Add two new sheets: Ambassis, Sheet3
# ---------------------------------------------------------------
$ws3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws2)
$ws3.Name = "Ambassis"

$ws4 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws3)
$ws4.Name = "Sheet3"

# ---------------------------------------------------------------
# 3. Sheet1 "Spratelloides" -- insert new row 3, renumber data
# ---------------------------------------------------------------
$ws1.Rows("3:3").Insert()

$ws1.Range("A2").Value = "Albatross_delicatulus"

$ws1.Range("A3").Value = "Albatross_gracilis"
$ws1.Range("B3").Value = 0.0059858000000000003
$ws1.Range("C3").Formula = "=LOG10(B3)"
$ws1.Range("D3").Value = 3.1470501999999998

$ws1.Range("A2").Select()

# ---------------------------------------------------------------
# 4. Sheet2 "Gerres" -- update formula sharing & view
# ---------------------------------------------------------------
$ws2.Range("C4:C30").Formula = "=LOG10(B4)"
$ws2.Range("C31:C43").Formula = "=LOG10(B31)"
$ws2.Range("A1:D1").Select()

# ---------------------------------------------------------------
# 5. Sheet3 "Ambassis" -- populate with new data
# ---------------------------------------------------------------
$ws3.Range("A1").Value = "Species"
$ws3.Range("B1").Value = "a"
$ws3.Range("C1").Value = "log10a"
$ws3.Range("D1").Value = "b"

$ws3.Range("A2").Value = "Albatross"
$ws3.Range("B2").Value = 0.031066
$ws3.Range("C2").Formula = "=LOG10(B2)"
$ws3.Range("D2").Value = 2.9220799999999998

$ws3.Range("A3").Value = "urotaenia"
$ws3.Range("B3").Value = 0.0214
$ws3.Range("C3").Formula = "=LOG10(B3)"
$ws3.Range("D3").Value = 2.653

$ws3.Range("A4").Value = "interrupta"
$ws3.Range("B4").Value = 0.015
$ws3.Range("D4").Value = 2.77

$ws3.Range("A5").Value = "interrupta"
$ws3.Range("B5").Value = 0.0328
$ws3.Range("D5").Value = 2.7930000000000001

$ws3.Range("A6").Value = "natalensis"
$ws3.Range("B6").Value = 0.0261
$ws3.Range("D6").Value = 2.964

$ws3.Range("A7").Value = "interrupta"
$ws3.Range("B7").Value = 0.0131
$ws3.Range("D7").Value = 2.984

$ws3.Range("A8").Value = "gymnocephalus"
$ws3.Range("B8").Value = 0.0154
$ws3.Range("D8").Value = 2.9870000000000001

$ws3.Range("A9").Value = "gymnocephalus"
$ws3.Range("B9").Value = 0.011
$ws3.Range("D9").Value = 3.0070000000000001

$ws3.Range("A10").Value = "gymnocephalus"
$ws3.Range("B10").Value = 0.0173
$ws3.Range("D10").Value = 3.0830000000000002

$ws3.Range("A11").Value = "urotaenia"
$ws3.Range("B11").Value = 0.014
$ws3.Range("D11").Value = 3.23

$ws3.Range("C4:C11").Formula = "=LOG10(B4)"

$ws3.Range("F7").Select()
$ws3.Activate()

# ---------------------------------------------------------------
# 6. Sheet4 "Sheet3" -- empty, leave as created
# ---------------------------------------------------------------
